$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# EPBDS-1587 add convertors for BigInteger, BigDecimal. tests
#
# Adds a new "Datatype TestBigTypes" table (rows 22-25) below the existing
# datatype tables on the first sheet, mirroring the layout already used for
# "Datatype TestType2" / "Datatype TestType3" (rows 15-19).
# ---------------------------------------------------------------------------

# ---- Cell values --------------------------------------------------------
# (Entered in this exact order so new shared-string entries get registered
#  in the same order the original workbook uses.)
$ws.Range("B22").Value = "Datatype TestBigTypes"

$ws.Range("B23").Value = "BigInteger"
$ws.Range("C23").Value = "bigIntVal"
$ws.Range("D23").Value = 2000000000

$ws.Range("C24").Value = "bigDecVal"
$ws.Range("B24").Value = "BigDecimal"
$ws.Range("D24").Value = 1115.3699999999999

$ws.Range("B25").Value = "BigInteger"
$ws.Range("C25").Value = "bigIntVal2"

# ---- Structure --------------------------------------------------------
# Merge before formatting so the merge doesn't strip the inner borders that
# the subsequent PasteSpecial sets (matches B15:D15 / J15:L15 which stay on
# one uniform style even though merged).
$null = $ws.Range("B22:D22").Merge()

# ---- Formatting -----------------------------------------------------------
# Header-style row (thin border all around + centered), same look as the
# "Datatype Test..." header rows (B3:D3, B15:D15).
$ws.Range("B15:D15").Copy()
$ws.Range("B22:D22").PasteSpecial(-4122)

$ws.Range("B15:D15").Copy()
$ws.Range("B23:D23").PasteSpecial(-4122)

$ws.Range("B15").Copy()
$ws.Range("C24").PasteSpecial(-4122)

$ws.Range("B15").Copy()
$ws.Range("C25").PasteSpecial(-4122)

# Plain bordered cells (same look as B4/D4).
$ws.Range("B4").Copy()
$ws.Range("B24").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("D24").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("B25").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# C17 switches from the plain border style to the centered/bordered style.
$ws.Range("B15").Copy()
$ws.Range("C17").PasteSpecial(-4122)

# ---- Selection (matches the state the workbook was saved in) --------------
$null = $ws.Range("C25").Select()
